$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest scrape.
# D-column cells whose new text parses as a plain number need a forced
# text format so Excel stores them as strings (matching source data
# which uses dotted/locale price formatting), not as numeric values.

$ws.Range('D2').Value = '29.006.75'
$ws.Range('E2').Value = '  -4.15%  '
$ws.Range('D3').Value = '1.962.70'
$ws.Range('E3').Value = '  -6.10%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.22%  '
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4991'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.98%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4210'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.78%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.19'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09219'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.099'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.01'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -6.75%  '
$ws.Range('D13').Value = '1.975.49'
$ws.Range('E13').Value = '  -2.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.877'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.95%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.443'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.26%  '
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('E17').Value = '  -4.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '91.38'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -10.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06664'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.24'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -8.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.005'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.955'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Value = '29.038.21'
$ws.Range('E23').Value = '  -3.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.02'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.282'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.41%  '
$ws.Range('D26').Value = '2.218.83'
$ws.Range('E26').Value = '  -2.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '156.72'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.59'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.205'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -9.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.265'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -9.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '126.81'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.13%  '
$ws.Range('E32').Value = '  -7.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09851'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.529'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -8.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.772'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.66%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.674'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02427'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.299'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.05%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.941'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -11.32%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06313'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.64%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6451'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.00%  '
$ws.Range('E42').Value = '  -8.82%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1994'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -9.81%  '
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6220'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -8.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.31'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.181'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.288'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.463'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.66%  '
$ws.Range('E50').Value = '  -4.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06910'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.84%  '
